# Edit "E suite.xlsx": add Jira id column, rename TCIDs, reword descriptions,
# tidy row heights / column widths on the "Test Cases" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Insert a new column B ("Jira id") before the existing Description column ---
$ws.Columns.Item(2).Insert()

# --- Header row ---
$ws.Cells.Item(1,1).Value = "TCID"
$ws.Cells.Item(1,2).Value = "Jira id"
$ws.Cells.Item(1,3).Value = "Description"
$ws.Cells.Item(1,4).Value = "Runmode"
$ws.Cells.Item(1,5).Value = "Results"

# --- Data rows: TCID, Jira id, Description, Runmode, Results ---
$ws.Cells.Item(2,1).Value = "TestCase_E1"
$ws.Cells.Item(2,2).Value = "OPQA-256"
$ws.Cells.Item(2,3).Value = "Verify that user is able to add a document to watchlist from search results page"
$ws.Cells.Item(2,4).Value = "Y"
$ws.Cells.Item(2,5).Value = "SKIP"

$ws.Cells.Item(3,1).Value = "TestCase_E2"
$ws.Cells.Item(3,2).Value = "OPQA-259"
$ws.Cells.Item(3,3).Value = "Verify that user is able to add a document to watchlist from document page"
$ws.Cells.Item(3,4).Value = "Y"
$ws.Cells.Item(3,5).Value = "SKIP"

$ws.Cells.Item(4,1).Value = "TestCase_E3"
$ws.Cells.Item(4,2).Value = "OPQA-260"
$ws.Cells.Item(4,3).Value = "Verify that user is able to delete a document from watchlist"
$ws.Cells.Item(4,4).Value = "Y"
$ws.Cells.Item(4,5).Value = "SKIP"

$ws.Cells.Item(5,1).Value = "TestCase_E4"
$ws.Cells.Item(5,2).Value = "OPQA-261"
$ws.Cells.Item(5,3).Value = "Verify that user is able to unwatch a document from search results page"
$ws.Cells.Item(5,4).Value = "Y"
$ws.Cells.Item(5,5).Value = "SKIP"

$ws.Cells.Item(6,1).Value = "TestCase_E5"
$ws.Cells.Item(6,2).Value = "OPQA-262"
$ws.Cells.Item(6,3).Value = "Verify that user is able to unwatch a document from document(Record View) page"
$ws.Cells.Item(6,4).Value = "Y"
$ws.Cells.Item(6,5).Value = "FAIL"

$ws.Cells.Item(7,1).Value = "TestCase_E6"
$ws.Cells.Item(7,2).Value = "OPQA-264"
$ws.Cells.Item(7,3).Value = "Verify that the following fields are getting displayed for each document in watchlist page:`na)Times cited`nb)Comments`nc)Views"
$ws.Cells.Item(7,4).Value = "Y"
$ws.Cells.Item(7,5).Value = "SKIP"

$ws.Cells.Item(8,1).Value = "TestCase_E7"
$ws.Cells.Item(8,2).Value = "OPQA-265"
$ws.Cells.Item(8,3).Value = "Verify that document count gets decreased in the watchlist page when a document is deleted from watchlist"
$ws.Cells.Item(8,4).Value = "Y"
$ws.Cells.Item(8,5).Value = "SKIP"

$ws.Cells.Item(9,1).Value = "TestCase_E8"
$ws.Cells.Item(9,2).Value = "OPQA-267"
$ws.Cells.Item(9,3).Value = "Verify that MORE button doesn't get displayed if number of documents in watchlist page is less than or equal to 10"
$ws.Cells.Item(9,4).Value = "Y"
$ws.Cells.Item(9,5).Value = "SKIP"

$ws.Cells.Item(10,1).Value = "TestCase_E9"
$ws.Cells.Item(10,2).Value = "OPQA-268"
$ws.Cells.Item(10,3).Value = "Verify that MORE button is present in watchlist page if total search results is more than 10`nVerify that MORE button is working correctly in watchlist page`n"
$ws.Cells.Item(10,4).Value = "Y"
$ws.Cells.Item(10,5).Value = "SKIP"

$ws.Cells.Item(11,1).Value = "TestCase_E10"
$ws.Cells.Item(11,2).Value = "OPQA-269"
$ws.Cells.Item(11,3).Value = "Verify that app navigates to correct page when user navigates back from document page"
$ws.Cells.Item(11,4).Value = "Y"
$ws.Cells.Item(11,5).Value = "SKIP"

# --- Row heights: row 9 loses its old custom height, row 10 becomes 45 ---
$ws.Rows.Item(9).AutoFit()
$ws.Rows.Item(10).RowHeight = 45

# --- Column widths (closest values the engine's quantized ColumnWidth can hit) ---
$ws.Columns.Item(1).ColumnWidth = 15.665
$ws.Columns.Item(2).ColumnWidth = 23.71
$ws.Columns.Item(3).ColumnWidth = 104.5
$ws.Columns.Item(4).ColumnWidth = 24.5
$ws.Columns.Item(5).ColumnWidth = 6.5

# --- Selection ---
$ws.Range("D2:D11").Select()
